$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "25.803.21"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.635.45"
$ws.Range("E3").Value = "  -0.14%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.16%  "

# Row 5 - BNB
Set-TextValue "D5" "215.45"
$ws.Range("E5").Value = "  +0.08%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.68%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.12%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.11%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.91%  "

# Row 10 - Solana
Set-TextValue "D10" "19.66"
$ws.Range("E10").Value = "  -3.55%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0791"
$ws.Range("E11").Value = "  +1.49%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  -0.06%  "

# Row 13 - was WrappedliquidstakedEther2.0, now WrappedEther (rows 13/14 swapped content)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.635.38"
$ws.Range("E13").Value = "  -0.45%  "

# Row 14 - was WrappedEther, now WrappedliquidstakedEther2.0
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D14" "1.860.09"
$ws.Range("E14").Value = "  -0.21%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.04%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0₃0769"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17 - Litecoin
Set-TextValue "D17" "62.89"

# Row 18 - WrappedBTC
Set-TextValue "D18" "25.815.73"
$ws.Range("E18").Value = "  -0.32%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.15%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +1.63%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "193.91"
$ws.Range("E21").Value = "  -0.20%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -0.02%  "

# Row 23 - Chainlink
Set-TextValue "D23" "6.30"
$ws.Range("E23").Value = "  +2.43%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  -0.08%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +1.94%  "

# Row 26 - Monero
Set-TextValue "D26" "142.56"
$ws.Range("E26").Value = "  +3.18%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +1.06%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +0.96%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +0.05%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.25%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.61%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +1.75%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -0.19%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +0.44%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +0.04%  "

# Row 37 - Maker
Set-TextValue "D37" "1.136.84"
$ws.Range("E37").Value = "  -0.07%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  -1.70%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -1.87%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.60%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.31%  "

# Row 42 - FraxShare
Set-TextValue "D42" "5.58"
$ws.Range("E42").Value = "  +1.96%  "

# Row 43 - Quant
Set-TextValue "D43" "100.53"
$ws.Range("E43").Value = "  +0.95%  "

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  +0.65%  "

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.770.48"

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  -4.13%  "

# Row 47 - Aave
Set-TextValue "D47" "55.26"
$ws.Range("E47").Value = "  -0.66%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  -0.29%  "

# Row 49 - Mantle
Set-TextValue "D49" "0.416"
$ws.Range("E49").Value = "  -2.08%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "7.54"
$ws.Range("E50").Value = "  -3.14%  "
